$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet (workbook.xml <sheet name=.../>)
$ws.Name = "Reroute To Sea"

# 2. Update the scenario title in B1
$ws.Range("B1").Value = "Reroute Cargo from HKG - DAL to HKG - SEA"

# 3. Remove the blank spacer row (old row 16) so the "When/Assert/Returns" block
#    shifts up by one (old rows 18,20,21,22,23 -> new rows 17,19,20,21,22)
$ws.Rows(16).Delete() | Out-Null

# 4. "Returns" row (now row 20): was "Rerouted Cargo" / "Cargo"
$ws.Range("B20").Value = "Returns"
$ws.Range("B20").Style = "40% - Accent4"

# 5. Origin row (now row 21): "Origin of" -> "Origin", insert "=" then shift HKG right
$ws.Range("E21").Value2 = $ws.Range("D21").Value2
$ws.Range("E21").Style = $ws.Range("D21").Style
$ws.Range("C21").Value = "Origin"
$ws.Range("C21").Style = "40% - Accent6"
$ws.Range("D21").Value = "'="
$ws.Range("D21").HorizontalAlignment = -4108

# 6. Destination row (now row 22): "Destination of" -> "Destination", insert "=" then shift SEA right
$ws.Range("E22").Value2 = $ws.Range("D22").Value2
$ws.Range("E22").Style = $ws.Range("D22").Style
$ws.Range("C22").Value = "Destination"
$ws.Range("C22").Style = "40% - Accent6"
$ws.Range("D22").Value = "'="
$ws.Range("D22").HorizontalAlignment = -4108

# 7. Restore the selection to B4 (was B2 before the edit)
$ws.Range("B4").Select() | Out-Null
